$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.418.07"
$ws.Range("E2").Value = "  -1.26%  "

$ws.Range("D3").Value = "2.528.94"
$ws.Range("E3").Value = "  -1.25%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'308.69"
$ws.Range("E5").Value = "  -2.26%  "

$ws.Range("D6").Value = "'99.18"
$ws.Range("E6").Value = "  +2.72%  "

$ws.Range("E7").Value = "  -1.35%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -1.95%  "

$ws.Range("D10").Value = "'35.82"
$ws.Range("E10").Value = "  +0.79%  "

$ws.Range("E11").Value = "  -1.32%  "

$ws.Range("D12").Value = "'7.32"
$ws.Range("E12").Value = "  -1.56%  "

$ws.Range("D14").Value = "2.917.59"
$ws.Range("E14").Value = "  -1.33%  "

$ws.Range("D15").Value = "'15.69"
$ws.Range("E15").Value = "  +4.05%  "

$ws.Range("D16").Value = "2.601.20"
$ws.Range("E16").Value = "  +2.41%  "

$ws.Range("D17").Value = "'0.819"
$ws.Range("E17").Value = "  -2.77%  "

$ws.Range("D18").Value = "42.423.26"
$ws.Range("E18").Value = "  -1.43%  "

$ws.Range("D19").Value = "'6.80"
$ws.Range("E19").Value = "  -0.42%  "

$ws.Range("D20").Value = "0.0₃0954"
$ws.Range("E20").Value = "  -0.62%  "

$ws.Range("D21").Value = "'12.16"
$ws.Range("E21").Value = "  -2.94%  "

$ws.Range("D22").Value = "'69.10"
$ws.Range("E22").Value = "  -0.39%  "

$ws.Range("D23").Value = "'244.17"
$ws.Range("E23").Value = "  -3.27%  "

$ws.Range("E24").Value = "  -2.26%  "

$ws.Range("E25").Value = "  -1.00%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").Value = "'25.90"
$ws.Range("E27").Value = "  -3.34%  "

$ws.Range("D28").Value = "'2.30"
$ws.Range("E28").Value = "  -5.85%  "

$ws.Range("D29").Value = "'39.22"
$ws.Range("E29").Value = "  -1.56%  "

$ws.Range("D30").Value = "'10.10"
$ws.Range("E30").Value = "  -0.94%  "

$ws.Range("D31").Value = "'156.21"
$ws.Range("E31").Value = "  +1.14%  "

$ws.Range("D32").Value = "'5.70"
$ws.Range("E32").Value = "  -1.90%  "

$ws.Range("E33").Value = "  +14.29%  "

$ws.Range("D34").Value = "'0.0792"
$ws.Range("E34").Value = "  -1.54%  "

$ws.Range("E35").Value = "  -3.16%  "

$ws.Range("D36").Value = "'2.02"
$ws.Range("E36").Value = "  -5.05%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'3.17"
$ws.Range("E37").Value = "  -7.02%  "

$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").Value = "'18.12"
$ws.Range("E38").Value = "  -4.99%  "

$ws.Range("E39").Value = "  -0.57%  "

$ws.Range("E40").Value = "  +0.08%  "

$ws.Range("E41").Value = "  +9.50%  "

$ws.Range("D42").Value = "'22.02"
$ws.Range("E42").Value = "  -1.70%  "

$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("D44").Value = "'3.30"
$ws.Range("E44").Value = "  +1.30%  "

$ws.Range("E45").Value = "  -2.29%  "

$ws.Range("D46").Value = "1.966.51"
$ws.Range("E46").Value = "  -1.75%  "

$ws.Range("E47").Value = "  -1.15%  "

$ws.Range("D48").Value = "2.771.33"
$ws.Range("E48").Value = "  -1.46%  "

$ws.Range("D49").Value = "'80.96"
$ws.Range("E49").Value = "  -2.74%  "

$ws.Range("B50").Value = "SEI"
$ws.Range("C50").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D50").Value = "'0.859"
$ws.Range("E50").Value = "  +11.02%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.192"
$ws.Range("E51").Value = "  -0.44%  "
